$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("V2").Value = 1.63

# Row 10
$ws.Range("H10").Value = 3.3
$ws.Range("J10").Value = 2.37
$ws.Range("S10").Value = 1.53
$ws.Range("T10").Value = 2.38
$ws.Range("AH10").Value = 12
$ws.Range("AI10").Value = 29
$ws.Range("AJ10").Value = 21
$ws.Range("AO10").Value = 9
$ws.Range("AT10").Value = 2.37
$ws.Range("AZ10").Value = 151

# Row 11
$ws.Range("G11").Value = 1.38
$ws.Range("S11").Value = 1.4

# Row 12
$ws.Range("G12").Value = 1.75
$ws.Range("I12").Value = 4.33
$ws.Range("J12").Value = 2.3
$ws.Range("U12").Value = 1.5
$ws.Range("V12").Value = 2.5
$ws.Range("AA12").Value = 13
$ws.Range("AD12").Value = 7.5
$ws.Range("AK12").Value = 41
$ws.Range("AO12").Value = 9

# Row 13
$ws.Range("N13").Value = 5.45
